$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full target dataset for rows 2-31 (the workbook was rerun with new/extended
# count data): columns are Row, A (index), B (timestamp), C, D, E, F, G, H.
$data = @(
    @(2, 0, 45271.33333333334, 1, 78, 1, 37, 2, 31),
    @(3, 1, 45271.34027777778, 0, 76, 2, 27, 2, 21),
    @(4, 2, 45271.34722222222, 1, 80, 2, 29, 2, 14),
    @(5, 3, 45271.35416666666, 0, 88, 1, 20, 1, 15),
    @(6, 4, 45271.36111111111, 1, 67, 1, 33, 1, 21),
    @(7, 5, 45271.36805555555, 0, 63, 3, 21, 1, 25),
    @(8, 6, 45377.66666666666, 2, 63, 0, 25, 2, 16),
    @(9, 7, 45377.67361111111, 2, 48, 0, 16, 1, 23),
    @(10, 8, 45377.68055555555, 1, 60, 0, 20, 1, 22),
    @(11, 9, 45377.6875, 2, 67, 0, 26, 1, 24),
    @(12, 10, 45377.69444444445, 3, 82, 0, 13, 2, 24),
    @(13, 11, 45377.70138888889, 2, 63, 0, 19, 3, 26),
    @(14, 12, 45391.33333333334, 3, 78, 0, 41, 0, 25),
    @(15, 13, 45391.34027777778, 0, 55, 1, 31, 2, 36),
    @(16, 14, 45391.34722222222, 3, 76, 0, 21, 1, 35),
    @(17, 15, 45391.35416666666, 1, 78, 0, 31, 1, 21),
    @(18, 16, 45391.36111111111, 1, 53, 3, 30, 2, 27),
    @(19, 17, 45391.36805555555, 4, 59, 1, 26, 3, 24),
    @(20, 18, 45391.66666666666, 1, 103, 2, 21, 1, 16),
    @(21, 19, 45391.67361111111, 2, 84, 1, 19, 1, 16),
    @(22, 20, 45391.68055555555, 2, 70, 3, 21, 1, 19),
    @(23, 21, 45391.6875, 1, 98, 1, 44, 2, 13),
    @(24, 22, 45391.69444444445, 1, 73, 3, 17, 0, 16),
    @(25, 23, 45391.70138888889, 4, 88, 0, 19, 3, 12),
    @(26, 24, 45391.83333333334, 1, 69, 0, 14, 2, 9),
    @(27, 25, 45391.84027777778, 0, 35, 0, 10, 1, 11),
    @(28, 26, 45391.84722222222, 1, 46, 2, 5, 0, 9),
    @(29, 27, 45391.85416666666, 0, 39, 0, 8, 0, 6),
    @(30, 28, 45391.86111111111, 1, 43, 1, 11, 3, 7),
    @(31, 29, 45391.86805555555, 0, 56, 0, 9, 0, 10)
)

# Rows 14-31 do not exist yet in the source sheet. Copy the formatting
# (number formats, borders, bold, alignment) from the last existing data
# row (13) down to row 31 before writing the new values, so the new rows
# match the existing style pattern (A: style 1, B: style 2 date format).
$ws.Range("A13:H13").Copy($ws.Range("A14:H31"))

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 5).Value = $row[5]
    $ws.Cells.Item($r, 6).Value = $row[6]
    $ws.Cells.Item($r, 7).Value = $row[7]
    $ws.Cells.Item($r, 8).Value = $row[8]
}
